$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.948267335492653
$ws.Range("D2").Value = 4.42801119425376
$ws.Range("E2").Value = 12.36580534531365
$ws.Range("F2").Value = 22.94271677801676
$ws.Range("G2").Value = 3.616670821718984
$ws.Range("I2").Value = 19.27961602937421
$ws.Range("K2").Value = 9.783469811797227
$ws.Range("M2").Value = 13.7021872114221
$ws.Range("N2").Value = 17.94595292912218
$ws.Range("O2").Value = 20.42392332626658
$ws.Range("B3").Value = 5.837326059220213
$ws.Range("D3").Value = 4.381651027853876
$ws.Range("E3").Value = 12.21535898797659
$ws.Range("F3").Value = 22.8990112739512
$ws.Range("G3").Value = 3.618601388786172
$ws.Range("I3").Value = 19.37349864929071
$ws.Range("K3").Value = 9.348140205501846
$ws.Range("M3").Value = 13.46812139071586
$ws.Range("N3").Value = 18.00322132273351
$ws.Range("O3").Value = 20.44937688565315
$ws.Range("B4").Value = 5.76907114511098
$ws.Range("D4").Value = 4.352493403077984
$ws.Range("E4").Value = 12.12654945993534
$ws.Range("F4").Value = 22.87897147127876
$ws.Range("G4").Value = 3.619849952815116
$ws.Range("I4").Value = 19.43420215985683
$ws.Range("K4").Value = 9.068203445303821
$ws.Range("M4").Value = 13.32562517399524
$ws.Range("N4").Value = 18.04007561289459
$ws.Range("O4").Value = 20.47017988601452
$ws.Range("B5").Value = 5.741260187519431
$ws.Range("D5").Value = 4.340442815701769
$ws.Range("E5").Value = 12.09129994339731
$ws.Range("F5").Value = 22.87251950943274
$ws.Range("G5").Value = 3.620374692561144
$ws.Range("I5").Value = 19.45971023999928
$ws.Range("K5").Value = 8.951034769721874
$ws.Range("M5").Value = 13.26794260146882
$ws.Range("N5").Value = 18.05552055076792
$ws.Range("O5").Value = 20.47995570002183
$ws.Range("B6").Value = 5.736643528332349
$ws.Range("D6").Value = 4.338431802390224
$ws.Range("E6").Value = 12.08550488416095
$ws.Range("F6").Value = 22.87155184073621
$ws.Range("G6").Value = 3.620462789408518
$ws.Range("I6").Value = 19.46399244678339
$ws.Range("K6").Value = 8.93139505660468
$ws.Range("M6").Value = 13.25839005706573
$ws.Range("N6").Value = 18.05811096774686
$ws.Range("O6").Value = 20.48165729774984
$ws.Range("B7").Value = 5.768696012237593
$ws.Range("D7").Value = 4.352331560089262
$ws.Range("E7").Value = 12.12607020502815
$ws.Range("F7").Value = 22.87887751022323
$ws.Range("G7").Value = 3.619856965031929
$ws.Range("I7").Value = 19.43454304718223
$ws.Range("K7").Value = 9.066635662582083
$ws.Range("M7").Value = 13.32484558052552
$ws.Range("N7").Value = 18.04028218015156
$ws.Range("O7").Value = 20.47030647282988
$ws.Range("B8").Value = 5.910065547528548
$ws.Range("D8").Value = 4.412173136818782
$ws.Range("E8").Value = 12.31321895538206
$ws.Range("F8").Value = 22.92624063307348
$ws.Range("G8").Value = 3.61732339708083
$ws.Range("I8").Value = 19.31135281140978
$ws.Range("K8").Value = 9.636040084721143
$ws.Range("M8").Value = 13.62127502455996
$ws.Range("N8").Value = 17.96534879380149
$ws.Range("O8").Value = 20.43162428340154
$ws.Range("B9").Value = 6.184592473256418
$ws.Range("D9").Value = 4.523796255867425
$ws.Range("E9").Value = 12.70643943349114
$ws.Range("F9").Value = 23.07270803207674
$ws.Range("G9").Value = 3.612854137188116
$ws.Range("I9").Value = 19.093974766686
$ws.Range("K9").Value = 10.64916979409083
$ws.Range("M9").Value = 14.20882730023885
$ws.Range("N9").Value = 17.83176710861391
$ws.Range("O9").Value = 20.39692495600039
$ws.Range("B10").Value = 6.382562217798113
$ws.Range("D10").Value = 4.602013694741309
$ws.Range("E10").Value = 13.00844453489782
$ws.Range("F10").Value = 23.21242157323974
$ws.Range("G10").Value = 3.609871579216561
$ws.Range("I10").Value = 18.94891167636167
$ws.Range("K10").Value = 11.32700994188589
$ws.Range("M10").Value = 14.63980018546209
$ws.Range("N10").Value = 17.74169001198439
$ws.Range("O10").Value = 20.39662458055399
$ws.Range("B11").Value = 6.471403168377099
$ws.Range("D11").Value = 4.636712319607071
$ws.Range("E11").Value = 13.14804878080704
$ws.Range("F11").Value = 23.28279315061996
$ws.Range("G11").Value = 3.608579409638112
$ws.Range("I11").Value = 18.88607613352863
$ws.Range("K11").Value = 11.62040391248958
$ws.Range("M11").Value = 14.83479071875039
$ws.Range("N11").Value = 17.70244575917745
$ws.Range("O11").Value = 20.40196806094457
$ws.Range("B12").Value = 6.504837937601839
$ws.Range("D12").Value = 4.649719853891439
$ws.Range("E12").Value = 13.20117944286636
$ws.Range("F12").Value = 23.31040484507721
$ws.Range("G12").Value = 3.608099337481424
$ws.Range("I12").Value = 18.86273392520587
$ws.Range("K12").Value = 11.72931912684763
$ws.Range("M12").Value = 14.90840305049872
$ws.Range("N12").Value = 17.68783281484702
$ws.Range("O12").Value = 20.40477919750772
$ws.Range("B13").Value = 6.497646881686099
$ws.Range("D13").Value = 4.646924402822247
$ws.Range("E13").Value = 13.18972584493375
$ws.Range("F13").Value = 23.30441560145131
$ws.Range("G13").Value = 3.608202319199118
$ws.Range("I13").Value = 18.8677409947293
$ws.Range("K13").Value = 11.70596003960136
$ws.Range("M13").Value = 14.89256052450778
$ws.Range("N13").Value = 17.69096896120004
$ws.Range("O13").Value = 20.40413874956223
$ws.Range("B14").Value = 6.474158202273034
$ws.Range("D14").Value = 4.637785133152412
$ws.Range("E14").Value = 13.15241485663048
$ws.Range("F14").Value = 23.28504557188273
$ws.Range("G14").Value = 3.608539728806325
$ws.Range("I14").Value = 18.88414669989976
$ws.Range("K14").Value = 11.6294084386177
$ws.Range("M14").Value = 14.84085179691615
$ws.Range("N14").Value = 17.7012385793445
$ws.Range("O14").Value = 20.40218355303417
$ws.Range("B15").Value = 6.459742770606619
$ws.Range("D15").Value = 4.632169716409781
$ws.Range("E15").Value = 13.12959379391134
$ws.Range("F15").Value = 23.27330582717219
$ws.Range("G15").Value = 3.608747604453305
$ws.Range("I15").Value = 18.89425452453221
$ws.Range("K15").Value = 11.58223265243727
$ws.Range("M15").Value = 14.80914709687202
$ws.Range("N15").Value = 17.70756128193008
$ws.Range("O15").Value = 20.4010884955233
$ws.Range("B16").Value = 6.376728959145013
$ws.Range("D16").Value = 4.599727850287882
$ws.Range("E16").Value = 12.9993611512772
$ws.Range("F16").Value = 23.20795839205972
$ws.Range("G16").Value = 3.609957321687258
$ws.Range("I16").Value = 18.95308148520201
$ws.Range("K16").Value = 11.30753174456031
$ws.Range("M16").Value = 14.62702936416538
$ws.Range("N16").Value = 17.74428947993001
$ws.Range("O16").Value = 20.39638563003846
$ws.Range("B17").Value = 6.325467446037345
$ws.Range("D17").Value = 4.579595990241964
$ws.Range("E17").Value = 12.91999806627083
$ws.Range("F17").Value = 23.16960419190636
$ws.Range("G17").Value = 3.610715959207468
$ws.Range("I17").Value = 18.98997686505754
$ws.Range("K17").Value = 11.13515277951774
$ws.Range("M17").Value = 14.51498076068052
$ws.Range("N17").Value = 17.7672639318997
$ws.Range("O17").Value = 20.39490411936579
$ws.Range("B18").Value = 6.295870506041201
$ws.Range("D18").Value = 4.567933872038886
$ws.Range("E18").Value = 12.87456274564441
$ws.Range("F18").Value = 23.14818644788525
$ws.Range("G18").Value = 3.611158391697471
$ws.Range("I18").Value = 19.01149509611059
$ws.Range("K18").Value = 11.03459862305342
$ws.Range("M18").Value = 14.45043770359476
$ws.Range("N18").Value = 17.78064134314857
$ws.Range("O18").Value = 20.39456787167888
$ws.Range("B19").Value = 6.285831157259847
$ws.Range("D19").Value = 4.563971216293273
$ws.Range("E19").Value = 12.8592172194155
$ws.Range("F19").Value = 23.14104560591021
$ws.Range("G19").Value = 3.611309238137685
$ws.Range("I19").Value = 19.01883184834316
$ws.Range("K19").Value = 11.00031237065053
$ws.Range("M19").Value = 14.42857040230452
$ws.Range("N19").Value = 17.78519875336805
$ws.Range("O19").Value = 20.3945426267072
$ws.Range("B20").Value = 6.330936223710125
$ws.Range("D20").Value = 4.58174766874063
$ws.Range("E20").Value = 12.92842482790214
$ws.Range("F20").Value = 23.17362067175164
$ws.Range("G20").Value = 3.610634571589291
$ws.Range("I20").Value = 18.9860185619683
$ws.Range("K20").Value = 11.1536486442584
$ws.Range("M20").Value = 14.52691897107474
$ws.Range("N20").Value = 17.76480138853577
$ws.Range("O20").Value = 20.3950084377379
$ws.Range("B21").Value = 6.481063274207514
$ws.Range("D21").Value = 4.640473182861142
$ws.Range("E21").Value = 13.16336721755073
$ws.Range("F21").Value = 23.29070900986621
$ws.Range("G21").Value = 3.608440372904171
$ws.Range("I21").Value = 18.87931568333416
$ws.Range("K21").Value = 11.65195309784314
$ws.Range("M21").Value = 14.85604661844052
$ws.Range("N21").Value = 17.69821542150839
$ws.Range("O21").Value = 20.40273647046624
$ws.Range("B22").Value = 6.577957526235434
$ws.Range("D22").Value = 4.678081570278289
$ws.Range("E22").Value = 13.31843897467764
$ws.Range("F22").Value = 23.3728402112031
$ws.Range("G22").Value = 3.607060199418664
$ws.Range("I22").Value = 18.81221462625808
$ws.Range("K22").Value = 11.96486592110923
$ws.Range("M22").Value = 15.06979935133841
$ws.Range("N22").Value = 17.65614269189029
$ws.Range("O22").Value = 20.41237773391557
$ws.Range("B23").Value = 6.526365335254153
$ws.Range("D23").Value = 4.658081565216973
$ws.Range("E23").Value = 13.2355525189938
$ws.Range("F23").Value = 23.32849807435793
$ws.Range("G23").Value = 3.607791910879255
$ws.Range("I23").Value = 18.84778701820612
$ws.Range("K23").Value = 11.79903593310512
$ws.Range("M23").Value = 14.95586264497841
$ws.Range("N23").Value = 17.67846583612985
$ws.Range("O23").Value = 20.40681226429503
$ws.Range("B24").Value = 6.328464179350712
$ws.Range("D24").Value = 4.580775168730222
$ws.Range("E24").Value = 12.92461448814959
$ws.Range("F24").Value = 23.17180284945842
$ws.Range("G24").Value = 3.610671347360737
$ws.Range("I24").Value = 18.9878071556903
$ws.Range("K24").Value = 11.14529118124353
$ws.Range("M24").Value = 14.52152208928561
$ws.Range("N24").Value = 17.76591417757723
$ws.Range("O24").Value = 20.3949596696228
$ws.Range("B25").Value = 6.110834687952604
$ws.Range("D25").Value = 4.494243883778375
$ws.Range("E25").Value = 12.59754841813815
$ws.Range("F25").Value = 23.02739866717607
$ws.Range("G25").Value = 3.614010098192242
$ws.Range("I25").Value = 19.15020192299665
$ws.Range("K25").Value = 10.38656180348286
$ws.Range("M25").Value = 14.04968220159312
$ws.Range("N25").Value = 17.86648234540651
$ws.Range("O25").Value = 20.4018929221796
